# Apply the daily cryptos price/volume update.
# Values are written with a leading apostrophe to force Excel to
# keep them as plain text (matching the original inlineStr cells),
# then ClearFormats() strips the resulting quote-prefix style so no
# stray cell style ("s" attribute) is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.003.49"
$ws.Range("E2").Value = "'  -1.58%  "
$ws.Range("D3").Value = "'2.986.02"
$ws.Range("E3").Value = "'  +0.38%  "
$ws.Range("E4").Value = "'  +0.17%  "
$ws.Range("D5").Value = "'503.08"
$ws.Range("E5").Value = "'  +0.64%  "
$ws.Range("D6").Value = "'138.15"
$ws.Range("E6").Value = "'  +0.36%  "
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E8").Value = "'  -0.73%  "
$ws.Range("D9").Value = "'7.12"
$ws.Range("E9").Value = "'  -2.64%  "
$ws.Range("E10").Value = "'  -1.36%  "
$ws.Range("D11").Value = "'0.363"
$ws.Range("E11").Value = "'  +1.61%  "
$ws.Range("D12").Value = "'3.496.02"
$ws.Range("E12").Value = "'  +0.48%  "
$ws.Range("E13").Value = "'  -1.69%  "
$ws.Range("D14").Value = "'26.01"
$ws.Range("E14").Value = "'  +0.02%  "
$ws.Range("E15").Value = "'  +0.44%  "
$ws.Range("D16").Value = "'56.062.54"
$ws.Range("E16").Value = "'  -1.54%  "
$ws.Range("B17").Value = "'Polkadot"
$ws.Range("C17").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'6.01"
$ws.Range("E17").Value = "'  -0.55%  "
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'2.984.00"
$ws.Range("E18").Value = "'  +0.20%  "
$ws.Range("D19").Value = "'12.95"
$ws.Range("E19").Value = "'  +2.75%  "
$ws.Range("D20").Value = "'7.98"
$ws.Range("E20").Value = "'  +1.41%  "
$ws.Range("D21").Value = "'327.75"
$ws.Range("E21").Value = "'  +2.04%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("D23").Value = "'0.493"
$ws.Range("E23").Value = "'  +1.04%  "
$ws.Range("D24").Value = "'64.61"
$ws.Range("E24").Value = "'  +1.69%  "
$ws.Range("D25").Value = "'3.104.20"
$ws.Range("E25").Value = "'  +0.38%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "'  -0.28%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "'  -1.81%  "
$ws.Range("D28").Value = "'0.0₃0898"
$ws.Range("E28").Value = "'  +0.56%  "
$ws.Range("E29").Value = "'  -1.92%  "
$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "'  -1.09%  "
$ws.Range("E31").Value = "'  +0.37%  "
$ws.Range("D32").Value = "'20.21"
$ws.Range("E32").Value = "'  +0.03%  "
$ws.Range("E33").Value = "'  -0.70%  "
$ws.Range("D34").Value = "'153.68"
$ws.Range("E34").Value = "'  -0.95%  "
$ws.Range("E35").Value = "'  -1.83%  "
$ws.Range("D36").Value = "'5.72"
$ws.Range("E36").Value = "'  -1.17%  "
$ws.Range("D37").Value = "'25.27"
$ws.Range("E37").Value = "'  +4.66%  "
$ws.Range("E38").Value = "'  -0.86%  "
$ws.Range("D39").Value = "'0.0658"
$ws.Range("E39").Value = "'  -1.49%  "
$ws.Range("D40").Value = "'3.022.82"
$ws.Range("E40").Value = "'  +0.58%  "
$ws.Range("D41").Value = "'36.62"
$ws.Range("E41").Value = "'  -2.75%  "
$ws.Range("E42").Value = "'  +0.08%  "
$ws.Range("D43").Value = "'3.77"
$ws.Range("E43").Value = "'  +1.02%  "
$ws.Range("D44").Value = "'0.648"
$ws.Range("E44").Value = "'  +1.57%  "
$ws.Range("D45").Value = "'2.159.51"
$ws.Range("E45").Value = "'  -1.50%  "
$ws.Range("E46").Value = "'  -3.15%  "
$ws.Range("D47").Value = "'5.85"
$ws.Range("E47").Value = "'  -1.67%  "
$ws.Range("D48").Value = "'0.921"
$ws.Range("E48").Value = "'  -2.16%  "
$ws.Range("B49").Value = "'VeChain"
$ws.Range("C49").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0236"
$ws.Range("E49").Value = "'  +0.74%  "
$ws.Range("B50").Value = "'InjectiveProtocol"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'19.66"
$ws.Range("E50").Value = "'  +2.22%  "
$ws.Range("D51").Value = "'0.0849"
$ws.Range("E51").Value = "'  -3.44%  "

$ws.Range("B2:E51").ClearFormats()
